# Commit: "change to currency data point"
#
# The "Component" (column F) for the financial figures Equity, Debt,
# Balance sheet total and EVIC (rows 5-8) changes from "Number" to
# "Currency", turning them into currency data points. Since the
# component itself now represents a currency value, the previously
# separate "Unit" (column H) value of "Currency" is no longer needed
# and is cleared (its formatting/style is kept).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Framework Data Model")

$ws.Range("F5").Value = "Currency"
$ws.Range("F6").Value = "Currency"
$ws.Range("F7").Value = "Currency"
$ws.Range("F8").Value = "Currency"

$ws.Range("H5").Value = $null
$ws.Range("H6").Value = $null
$ws.Range("H7").Value = $null
$ws.Range("H8").Value = $null
